# Natmi following Dr Hou advice:
# the LR-pair analysis now includes a third "ECs" cluster alongside the
# existing "FAPs" and "sCs" clusters, so the sending/target cluster grid
# grows from 2x2 (rows 2-5) to 3x3 (rows 2-10) for the Bmp4/Rgmb pair.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full data table (headers in row 1 stay as-is; rows 2-10 get rewritten)
$data = @(
    ,("ECs", "Bmp4", "Rgmb", "ECs", 2, 0.6666666666666666, 5.448315, 16.344945, 0.1199618029178375, 0.1199618029178375, 2, 0.6666666666666666, 12.531022, 37.593066, 0.1764330600309793, 0.1764330600309793, 68.27295512793, 614.4565961513699, 0.02116522797562734, 0.02116522797562734)
    ,("ECs", "Bmp4", "Rgmb", "FAPs", 2, 0.6666666666666666, 5.448315, 16.344945, 0.1199618029178375, 0.1199618029178375, 3, 1, 37.18315166666667, 111.549455, 0.5235277082863637, 0.5235277082863639, 202.585522972775, 1823.269706754975, 0.06280332776347591, 0.06280332776347593)
    ,("ECs", "Bmp4", "Rgmb", "sCs", 2, 0.6666666666666666, 5.448315, 16.344945, 0.1199618029178375, 0.1199618029178375, 3, 1, 21.310055, 63.930165, 0.3000392316826569, 0.3000392316826569, 116.103892307325, 1044.935030765925, 0.03599324717873428, 0.03599324717873428)
    ,("FAPs", "Bmp4", "Rgmb", "ECs", 3, 1, 31.28251333333334, 93.84754000000001, 0.6887829905701046, 0.6887829905701045, 2, 0.6666666666666666, 12.531022, 37.593066, 0.1764330600309793, 0.1764330600309793, 392.0018627952934, 3528.01676515764, 0.1215240907235727, 0.1215240907235727)
    ,("FAPs", "Bmp4", "Rgmb", "FAPs", 3, 1, 31.28251333333334, 93.84754000000001, 0.6887829905701046, 0.6887829905701045, 3, 1, 37.18315166666667, 111.549455, 0.5235277082863637, 0.5235277082863639, 1163.182437787856, 10468.6419400907, 0.3605969805597949, 0.360596980559795)
    ,("FAPs", "Bmp4", "Rgmb", "sCs", 3, 1, 31.28251333333334, 93.84754000000001, 0.6887829905701046, 0.6887829905701045, 3, 1, 21.310055, 63.930165, 0.3000392316826569, 0.3000392316826569, 666.6320796715668, 5999.6887170441, 0.2066619192867369, 0.2066619192867369)
    ,("sCs", "Bmp4", "Rgmb", "ECs", 3, 1, 8.686253333333333, 26.05876, 0.1912552065120579, 0.1912552065120579, 2, 0.6666666666666666, 12.531022, 37.593066, 0.1764330600309793, 0.1764330600309793, 108.8476316175733, 979.62868455816, 0.03374374133177925, 0.03374374133177926)
    ,("sCs", "Bmp4", "Rgmb", "FAPs", 3, 1, 8.686253333333333, 26.05876, 0.1912552065120579, 0.1912552065120579, 3, 1, 37.18315166666667, 111.549455, 0.5235277082863637, 0.5235277082863639, 322.9822751084222, 2906.8404759758, 0.1001273999630929, 0.1001273999630929)
    ,("sCs", "Bmp4", "Rgmb", "sCs", 3, 1, 8.686253333333333, 26.05876, 0.1912552065120579, 0.1912552065120579, 3, 1, 21.310055, 63.930165, 0.3000392316826569, 0.3000392316826569, 185.1045362772667, 1665.9408264954, 0.05738406521718573, 0.05738406521718573)
)

$startRow = 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $rowData = $data[$i]
    $r = $startRow + $i
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($r, $c + 1).Value = $rowData[$c]
    }
}
